$d = $word.ActiveDocument
$tbl = $d.Tables(1)
$nl = [char]11

# Row 1 Col 1: 35 x 39 -> 70 x 38
$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "70 x 38" + $nl + "  3    8" + $nl + "  ----" + $nl + "7|    |" + $nl + "0|    |"

# Row 1 Col 2: 32 x 90 -> 46 x 53
$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "46 x 53" + $nl + "  5    3" + $nl + "  ----" + $nl + "4|    |" + $nl + "6|    |"

# Row 1 Col 3: 29 x 30 -> 99 x 19
$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "99 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "9|    |" + $nl + "9|    |"

# Row 2 Col 1: 30 x 79 -> 85 x 40
$cell = $tbl.Cell(2, 1)
$cell.Range.Text = "85 x 40" + $nl + "  4    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "5|    |"

# Row 2 Col 2: 38 x 47 -> 65 x 94
$cell = $tbl.Cell(2, 2)
$cell.Range.Text = "65 x 94" + $nl + "  9    4" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |"

# Row 2 Col 3: 71 x 44 -> 51 x 19
$cell = $tbl.Cell(2, 3)
$cell.Range.Text = "51 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "5|    |" + $nl + "1|    |"

# Row 3 Col 1: 64 x 95 -> 30 x 74
$cell = $tbl.Cell(3, 1)
$cell.Range.Text = "30 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "3|    |" + $nl + "0|    |"

# Row 3 Col 2: 39 x 16 -> 95 x 12
$cell = $tbl.Cell(3, 2)
$cell.Range.Text = "95 x 12" + $nl + "  1    2" + $nl + "  ----" + $nl + "9|    |" + $nl + "5|    |"

# Row 3 Col 3: 16 x 94 -> 88 x 79
$cell = $tbl.Cell(3, 3)
$cell.Range.Text = "88 x 79" + $nl + "  7    9" + $nl + "  ----" + $nl + "8|    |" + $nl + "8|    |"

# Row 4 Col 1: 17 x 21 -> 60 x 32
$cell = $tbl.Cell(4, 1)
$cell.Range.Text = "60 x 32" + $nl + "  3    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "0|    |"

# Row 4 Col 2: 97 x 39 -> 47 x 92
$cell = $tbl.Cell(4, 2)
$cell.Range.Text = "47 x 92" + $nl + "  9    2" + $nl + "  ----" + $nl + "4|    |" + $nl + "7|    |"

# Row 4 Col 3: 21 x 92 -> 56 x 90
$cell = $tbl.Cell(4, 3)
$cell.Range.Text = "56 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "5|    |" + $nl + "6|    |"

# Row 5 Col 1: 42 x 35 -> 80 x 86
$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "80 x 86" + $nl + "  8    6" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"

# Row 5 Col 2: 44 x 74 -> 44 x 84
$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "44 x 84" + $nl + "  8    4" + $nl + "  ----" + $nl + "4|    |" + $nl + "4|    |"

# Row 5 Col 3: 85 x 99 -> 74 x 34
$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "74 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "7|    |" + $nl + "4|    |"
